$d = $word.ActiveDocument

# Package wrapper used with Range.InsertXML() below - lets us hand Word a
# literal WordprocessingML fragment (so the freshly-split "Castro Valley..."
# run keeps an explicit xml:space="preserve", matching the rest of the
# letterhead block) without disturbing the paragraph/run formatting that
# Word already derived for us.
function New-WordPackageXml([string]$bodyFragment) {
    return '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyFragment + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. Update the letter date -------------------------------------------
foreach ($para in $d.Paragraphs) {
    $rng = $para.Range
    if ($rng.Text.TrimEnd([char]13, [char]7) -eq "September 19, 2025") {
        $rng.Text = "September 21, 2025"
        break
    }
}

# --- 2. Split "<street>, <city> <state> <zip>" onto its own two lines ----
foreach ($para in $d.Paragraphs) {
    $rng = $para.Range
    if ($rng.Text.TrimEnd([char]13, [char]7) -eq "19958 Lorena Circle, Castro Valley CA 94546") {
        # Shrinking the existing run in place keeps this paragraph's own
        # <w:p> identity (paraId/rsid/etc.) and its xml:space="preserve",
        # and the embedded "`r" mints a brand-new paragraph right after it
        # (inheriting the same pPr/rPr) for the second line.
        $rng.Text = "19958 Lorena Circle`rCastro Valley, CA 94546"

        $cityPara = $para.Next()
        $cityRange = $cityPara.Range
        $cityXml = New-WordPackageXml(
            '<w:p>' +
              '<w:pPr>' +
                '<w:autoSpaceDE w:val="0"/>' +
                '<w:autoSpaceDN w:val="0"/>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
                  '<w:sz w:val="22"/>' +
                  '<w:szCs w:val="22"/>' +
                '</w:rPr>' +
              '</w:pPr>' +
              '<w:r>' +
                '<w:rPr>' +
                  '<w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' +
                  '<w:sz w:val="22"/>' +
                  '<w:szCs w:val="22"/>' +
                '</w:rPr>' +
                '<w:t xml:space="preserve">Castro Valley, CA 94546</w:t>' +
              '</w:r>' +
            '</w:p>'
        )
        $cityRange.InsertXML($cityXml)
        break
    }
}

# --- 3. Drop the blank "No Spacing" paragraph that used to sit right
#        after "Board of Directors" in the letterhead/signature block ----
foreach ($para in $d.Paragraphs) {
    $rng = $para.Range
    if ($rng.Text.TrimEnd([char]13, [char]7) -eq "") {
        $prev = $para.Previous()
        if ($prev -ne $null -and
            $prev.Range.Text.TrimEnd([char]13, [char]7) -eq "Lorena Circle Homeowners Association Board of Directors") {
            $rng.Delete()
            break
        }
    }
}
